$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B7").Value = "['Others', 'Hyundai KONA 64 kWh']"
$ws.Range("C7").Value = "[0.4, 0.3]"
$ws.Range("D7").Value = "[0.7500000000000001, 0.7000000000000001]"
$ws.Range("E7").Value = "[18.420694444444447, 25.600000000000005]"
$ws.Range("F7").Value = 44.02069444444445

$ws.Range("B8").Value = "['Others', 'VW ID.3', 'Hyundai KONA 64 kWh']"
$ws.Range("C8").Value = "[0.15, 0.1, 0.3]"
$ws.Range("D8").Value = "[0.7000000000000001, 0.8500000000000002, 0.7000000000000001]"
$ws.Range("E8").Value = "[28.946805555555557, 43.500000000000014, 25.600000000000005]"
$ws.Range("F8").Value = 98.04680555555558

$ws.Range("B9").Value = "['Fiat 500 E', 'Tesla MODEL 3', 'Others']"
$ws.Range("C9").Value = "[0.4, 0.05, 0.1]"
$ws.Range("D9").Value = "[0.8500000000000002, 0.9500000000000003, 0.8500000000000002]"
$ws.Range("E9").Value = "[10.800000000000004, 45.000000000000014, 39.47291666666668]"
$ws.Range("F9").Value = 95.2729166666667

$ws.Range("B12").Value = "['VW ID.4', 'TESLA MODEL Y']"
$ws.Range("C12").Value = "[0.25, 0.25]"
$ws.Range("D12").Value = "[0.8500000000000002, 0.7000000000000001]"
$ws.Range("E12").Value = "[46.20000000000002, 33.75000000000001]"
$ws.Range("F12").Value = 79.95000000000002

$ws.Range("B13").Value = "['Others', 'Opel CORSA', 'Others', 'VW ID.3', 'Others']"
$ws.Range("C13").Value = "[0.4, 0.25, 0.15, 0.2, 0.3]"
$ws.Range("D13").Value = "[0.9500000000000003, 0.7000000000000001, 0.9000000000000002, 0.8000000000000002, 0.7500000000000001]"
$ws.Range("E13").Value = "[28.946805555555567, 20.250000000000004, 39.47291666666668, 34.800000000000004, 23.683750000000007]"
$ws.Range("F13").Value = 147.1534722222223

$ws.Range("B17").Value = "['Audi E-TRON', 'Others']"
$ws.Range("C17").Value = "[0.4, 0.25]"
$ws.Range("D17").Value = "[0.7000000000000001, 0.8000000000000002]"
$ws.Range("E17").Value = "[25.500000000000004, 28.946805555555564]"
$ws.Range("F17").Value = 54.44680555555557

$ws.Range("B18").Value = "['Others', 'Opel CORSA']"
$ws.Range("C18").Value = "[0.3, 0.3]"
$ws.Range("D18").Value = "[0.8500000000000002, 0.7500000000000001]"
$ws.Range("E18").Value = "[28.946805555555567, 20.250000000000007]"
$ws.Range("F18").Value = 49.19680555555557

$ws.Range("B19").Value = "[]"
$ws.Range("C19").Value = "[]"
$ws.Range("D19").Value = "[]"
$ws.Range("E19").Value = "[]"
$ws.Range("F19").Value = 0

$ws.Range("B20").Value = "['Others', 'VW ID.5']"
$ws.Range("C20").Value = "[0.05, 0.3]"
$ws.Range("D20").Value = "[0.7000000000000001, 0.7000000000000001]"
$ws.Range("E20").Value = "[34.20986111111111, 30.800000000000004]"
$ws.Range("F20").Value = 65.00986111111112

$ws.Range("B21").Value = "['Tesla MODEL 3']"
$ws.Range("C21").Value = "[0.3]"
$ws.Range("D21").Value = "[0.8500000000000002]"
$ws.Range("E21").Value = "[27.500000000000014]"
$ws.Range("F21").Value = 27.50000000000001

$ws.Range("B32").Value = "['MINI Cooper SE', 'Renault ZOE', 'Others', 'Fiat 500 E']"
$ws.Range("C32").Value = "[0.15, 0.25, 0.35, 0.25]"
$ws.Range("D32").Value = "[0.9000000000000002, 0.7500000000000001, 0.7500000000000001, 0.8500000000000002]"
$ws.Range("E32").Value = "[21.675000000000004, 26.000000000000007, 21.052222222222227, 14.400000000000006]"
$ws.Range("F32").Value = 83.12722222222224

$ws.Range("B33").Value = "['VW ID.3', 'Audi E-TRON']"
$ws.Range("C33").Value = "[0.2, 0.15]"
$ws.Range("D33").Value = "[0.65, 0.8000000000000002]"
$ws.Range("E33").Value = "[26.1, 55.250000000000014]"
$ws.Range("F33").Value = 81.35000000000002

$ws.Range("B34").Value = "['VW ID.5']"
$ws.Range("C34").Value = "[0.3]"
$ws.Range("D34").Value = "[0.7500000000000001]"
$ws.Range("E34").Value = "[34.65000000000001]"
$ws.Range("F34").Value = 34.65000000000001

$ws.Range("B35").Value = "['MINI Cooper SE']"
$ws.Range("C35").Value = "[0.25]"
$ws.Range("D35").Value = "[0.9000000000000002]"
$ws.Range("E35").Value = "[18.785000000000007]"
$ws.Range("F35").Value = 18.78500000000001

$ws.Range("B36").Value = "['SKODA ENYAQ 77kWh']"
$ws.Range("C36").Value = "[0.15]"
$ws.Range("D36").Value = "[0.65]"
$ws.Range("E36").Value = "[38.5]"
$ws.Range("F36").Value = 38.5

$ws.Range("B37").Value = "['VW ID.5', 'Renault ZOE']"
$ws.Range("C37").Value = "[0.2, 0.35]"
$ws.Range("D37").Value = "[0.7500000000000001, 0.7500000000000001]"
$ws.Range("E37").Value = "[42.35, 20.800000000000008]"
$ws.Range("F37").Value = 63.15000000000001

$ws.Range("B38").Value = "['Tesla MODEL 3', 'Tesla MODEL 3', 'Tesla MODEL 3']"
$ws.Range("C38").Value = "[0.25, 0.4999999999999999, 0.45]"
$ws.Range("D38").Value = "[0.8500000000000002, 0.9000000000000002, 0.65]"
$ws.Range("E38").Value = "[30.00000000000001, 20.000000000000018, 10.0]"
$ws.Range("F38").Value = 60.00000000000003

$ws.Range("B41").Value = "['Tesla MODEL 3', 'Others']"
$ws.Range("C41").Value = "[0.2, 0.45]"
$ws.Range("D41").Value = "[0.9500000000000003, 0.6]"
$ws.Range("E41").Value = "[37.500000000000014, 7.8945833333333315]"
$ws.Range("F41").Value = 45.39458333333334

$ws.Range("B42").Value = "['Dacia SPRING', 'Fiat 500 E']"
$ws.Range("C42").Value = "[0.15, 0.25]"
$ws.Range("D42").Value = "[0.7000000000000001, 0.9000000000000002]"
$ws.Range("E42").Value = "[14.740000000000002, 15.600000000000005]"
$ws.Range("F42").Value = 30.34000000000001

$ws.Range("B43").Value = "['VW ID.4', 'VW ID.3', 'Others']"
$ws.Range("C43").Value = "[0.1, 0.1, 0.1]"
$ws.Range("D43").Value = "[0.7000000000000001, 0.9500000000000003, 1.0]"
$ws.Range("E43").Value = "[46.20000000000001, 49.30000000000002, 47.3675]"
$ws.Range("F43").Value = 142.8675

$ws.Range("B44").Value = "[]"
$ws.Range("C44").Value = "[]"
$ws.Range("D44").Value = "[]"
$ws.Range("E44").Value = "[]"
$ws.Range("F44").Value = 0
